# Weekly refresh of the cryptocurrency price/volume table.
# Mirrors the upstream GitHub Actions job that re-scrapes coinranking.com
# and rewrites the "Price" / "Volume(1h)" columns (plus, this time, two rows
# that changed rank order: Maker<->Stacks and Cosmos<->ONDO).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some new "Price" strings look like plain decimals (e.g. "7.32", "1.00",
# "0.0235") and Excel would otherwise silently reinterpret them as numbers.
# Force those specific cells to Text format first so the literal string is kept,
# matching the inline-string cells already used throughout this column.
$textPriceCells = @(
    "D5", "D6", "D8", "D9", "D11", "D14", "D15", "D17", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D30", "D33", "D34", "D35", "D36", "D38", "D39", "D41", "D43", "D44", "D45", "D47", "D48", "D49", "D50"
)
foreach ($addr in $textPriceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "57.152.53"
$ws.Range("E2").Value = "  -1.11%  "
$ws.Range("D3").Value = "2.990.29"
$ws.Range("E3").Value = "  -2.06%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "501.45"
$ws.Range("E5").Value = "  -4.50%  "
$ws.Range("D6").Value = "138.46"
$ws.Range("E6").Value = "  -2.60%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "0.431"
$ws.Range("E8").Value = "  -3.37%  "
$ws.Range("D9").Value = "7.32"
$ws.Range("E9").Value = "  -4.04%  "
$ws.Range("E10").Value = "  -4.53%  "
$ws.Range("D11").Value = "0.360"
$ws.Range("E11").Value = "  -3.10%  "
$ws.Range("D12").Value = "3.507.84"
$ws.Range("E12").Value = "  -1.84%  "
$ws.Range("E13").Value = "  -2.37%  "
$ws.Range("D14").Value = "26.21"
$ws.Range("E14").Value = "  -2.73%  "
$ws.Range("D15").Value = "0.0000160"
$ws.Range("E15").Value = "  -6.13%  "
$ws.Range("D16").Value = "57.216.34"
$ws.Range("E16").Value = "  -1.01%  "
$ws.Range("D17").Value = "6.08"
$ws.Range("E17").Value = "  -2.87%  "
$ws.Range("D18").Value = "2.993.67"
$ws.Range("E18").Value = "  -1.97%  "
$ws.Range("D19").Value = "12.68"
$ws.Range("E19").Value = "  -3.05%  "
$ws.Range("D20").Value = "7.89"
$ws.Range("E20").Value = "  -3.75%  "
$ws.Range("D21").Value = "320.86"
$ws.Range("E21").Value = "  -5.69%  "
$ws.Range("D22").Value = "1.00"
$ws.Range("E22").Value = "  -0.13%  "
$ws.Range("D23").Value = "5.75"
$ws.Range("E23").Value = "  +1.19%  "
$ws.Range("D24").Value = "0.493"
$ws.Range("E24").Value = "  -1.83%  "
$ws.Range("D25").Value = "63.45"
$ws.Range("E25").Value = "  -2.34%  "
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("E27").Value = "  -5.55%  "
$ws.Range("D28").Value = "0.0₃0898"
$ws.Range("E28").Value = "  -7.57%  "
$ws.Range("E29").Value = "  -5.44%  "
$ws.Range("D30").Value = "7.14"
$ws.Range("E30").Value = "  -3.47%  "
$ws.Range("E31").Value = "  -4.01%  "
$ws.Range("E32").Value = "  -4.89%  "
$ws.Range("D33").Value = "20.28"
$ws.Range("E33").Value = "  -4.01%  "
$ws.Range("D34").Value = "154.92"
$ws.Range("E34").Value = "  -1.10%  "
$ws.Range("D35").Value = "4.58"
$ws.Range("E35").Value = "  -3.60%  "
$ws.Range("D36").Value = "5.79"
$ws.Range("E36").Value = "  -3.22%  "
$ws.Range("E37").Value = "  -5.89%  "
$ws.Range("D38").Value = "24.43"
$ws.Range("E38").Value = "  -5.50%  "
$ws.Range("D39").Value = "0.0664"
$ws.Range("E39").Value = "  -6.09%  "
$ws.Range("D40").Value = "3.026.74"
$ws.Range("E40").Value = "  -2.00%  "
$ws.Range("D41").Value = "37.78"
$ws.Range("E41").Value = "  +0.14%  "
$ws.Range("E42").Value = "  +0.04%  "
$ws.Range("D43").Value = "3.76"
$ws.Range("E43").Value = "  -3.40%  "
$ws.Range("D44").Value = "0.646"
$ws.Range("E44").Value = "  -2.49%  "
$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D45").Value = "1.39"
$ws.Range("E45").Value = "  -6.17%  "
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "2.196.99"
$ws.Range("B47").Value = "ONDO"
$ws.Range("C47").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D47").Value = "0.943"
$ws.Range("E47").Value = "  -8.51%  "
$ws.Range("B48").Value = "Cosmos"
$ws.Range("C48").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D48").Value = "5.96"
$ws.Range("E48").Value = "  -1.55%  "
$ws.Range("D49").Value = "0.0235"
$ws.Range("E49").Value = "  -4.75%  "
$ws.Range("D50").Value = "19.32"
$ws.Range("E50").Value = "  -4.50%  "
$ws.Range("E51").Value = "  -10.92%  "
